$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.647.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.728.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.07%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.26"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.46%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.548"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.726.81"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.34"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.157"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.68"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.225.71"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.636.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.725.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.91"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "374.78"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.65"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.51"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.98"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.57%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.60"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.867.88"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "590.05"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.11%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.30"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.44"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.74%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.44%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.92"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.95"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.380"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.90"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.47"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.98"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.66"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0312"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.603"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "155.43"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.93"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.77%  "
